# Advance the Dijkstra-step animation on Feuil1 by one more step:
#  - M15 ("0 11  13", orange) is cleared out (its value has been superseded)
#  - L17 gains the new orange "current" label  "1 14 7"
#  - M17's old orange label "1 11 13" turns green (finalised) and becomes "0 16 10"
#  - L18's old green "0 14 7" is cleared back to a plain cell
#  - M18 becomes the new green (finalised) label "1 16 10"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- M15: was orange "0 11  13" -> cleared to a plain white/bordered cell ---
$ws.Range("M15").Value = ""
$ws.Range("M15").Interior.Color = 16777215

# --- M17: was orange "1 11 13" -> green "0 16 10" (copy style from a green cell, then set text) ---
$ws.Range("C4").Copy($ws.Range("M17"))
$ws.Range("M17").Value = "0 16 10"

# --- M18: was empty (plain) -> green "1 16 10" ---
$ws.Range("C4").Copy($ws.Range("M18"))
$ws.Range("M18").Value = "1 16 10"

# --- L17: was empty (plain) -> orange "1 14 7" (copy style from an orange cell, then set text) ---
$ws.Range("C3").Copy($ws.Range("L17"))
$ws.Range("L17").Value = "1 14 7"

# --- L18: was green "0 14 7" -> cleared back to a plain bordered cell ---
$ws.Range("B15").Copy($ws.Range("L18"))
$ws.Range("L18").Value = ""

# --- Selection moved from D19 to L12 ---
$ws.Range("L12").Select()
